# Auto-generated edit script: Add data for 2023-08-26
# Updates YTD crime-count comparisons across Citywide Totals, By Neighborhood summary, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("B2").Value = 28
$ws.Range("D2").Value = 64
$ws.Range("E3").Value = 96
$ws.Range("F3").Value = 92
$ws.Range("I3").Value = 136
$ws.Range("B9").Value = 255
$ws.Range("C9").Value = 325
$ws.Range("D9").Value = 294
$ws.Range("E9").Value = 287
$ws.Range("F9").Value = 378
$ws.Range("G9").Value = 332
$ws.Range("B10").Value = 847
$ws.Range("C10").Value = 1009
$ws.Range("D10").Value = 1187
$ws.Range("E10").Value = 1471
$ws.Range("F10").Value = 1515
$ws.Range("G10").Value = 733
$ws.Range("H10").Value = 363
$ws.Range("I10").Value = 585
$ws.Range("J10").Value = 478
$ws.Range("B11").Value = 1198
$ws.Range("C11").Value = 1448
$ws.Range("D11").Value = 1646
$ws.Range("E11").Value = 1910
$ws.Range("F11").Value = 2054
$ws.Range("G11").Value = 1218
$ws.Range("H11").Value = 831
$ws.Range("I11").Value = 1188
$ws.Range("J11").Value = 992

# --- Chinatown ---
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("B7").Value = 4
$ws.Range("B9").Value = 12

# --- Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("D7").Value = 30
$ws.Range("F8").Value = 94
$ws.Range("D9").Value = 71
$ws.Range("F9").Value = 144

# --- Chatham ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("G6").Value = 12
$ws.Range("G8").Value = 30

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("B7").Value = 13
$ws.Range("D7").Value = 16
$ws.Range("F7").Value = 16
$ws.Range("G7").Value = 8
$ws.Range("F8").Value = 28
$ws.Range("B9").Value = 41
$ws.Range("D9").Value = 56
$ws.Range("F9").Value = 64
$ws.Range("G9").Value = 30

# --- Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 22
$ws.Range("B8").Value = 110
$ws.Range("C8").Value = 193
$ws.Range("E8").Value = 427
$ws.Range("F8").Value = 410
$ws.Range("H8").Value = 63
$ws.Range("B9").Value = 143
$ws.Range("C9").Value = 231
$ws.Range("E9").Value = 479
$ws.Range("F9").Value = 465
$ws.Range("H9").Value = 118
$ws.Range("I9").Value = 233

# --- Old Town ---
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 24

# --- Little Italy, UIC ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J5").Value = 6
$ws.Range("J7").Value = 21

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("B7").Value = 35
$ws.Range("F7").Value = 15
$ws.Range("H7").Value = 8
$ws.Range("C8").Value = 69
$ws.Range("E8").Value = 72
$ws.Range("F8").Value = 90
$ws.Range("G19").Value = 30
$ws.Range("B21").Value = 12
$ws.Range("D23").Value = 13
$ws.Range("C28").Value = 97
$ws.Range("F28").Value = 82
$ws.Range("C30").Value = 21
$ws.Range("D32").Value = 71
$ws.Range("F32").Value = 144
$ws.Range("B36").Value = 41
$ws.Range("D36").Value = 56
$ws.Range("F36").Value = 64
$ws.Range("G36").Value = 30
$ws.Range("B41").Value = 12
$ws.Range("E45").Value = 14
$ws.Range("J46").Value = 3
$ws.Range("F47").Value = 53
$ws.Range("H47").Value = 27
$ws.Range("J47").Value = 25
$ws.Range("G48").Value = 7
$ws.Range("J50").Value = 21
$ws.Range("B53").Value = 143
$ws.Range("C53").Value = 231
$ws.Range("E53").Value = 479
$ws.Range("F53").Value = 465
$ws.Range("H53").Value = 118
$ws.Range("I53").Value = 233
$ws.Range("I61").Value = 13
$ws.Range("J61").Value = 14
$ws.Range("D62").Value = 19
$ws.Range("F62").Value = 18
$ws.Range("F68").Value = 28
$ws.Range("F69").Value = 2
$ws.Range("I70").Value = 24
$ws.Range("B74").Value = 34
$ws.Range("D76").Value = 41
$ws.Range("E76").Value = 68
$ws.Range("F76").Value = 44
$ws.Range("C77").Value = 46
$ws.Range("E77").Value = 52
$ws.Range("J77").Value = 50
$ws.Range("D87").Value = 14
$ws.Range("E87").Value = 24
$ws.Range("B88").Value = 6
$ws.Range("C89").Value = 22
$ws.Range("I92").Value = 25
$ws.Range("B93").Value = 13
$ws.Range("F93").Value = 8
$ws.Range("B95").Value = 10
$ws.Range("B97").Value = 21
$ws.Range("B99").Value = 1198
$ws.Range("C99").Value = 1448
$ws.Range("D99").Value = 1646
$ws.Range("E99").Value = 1910
$ws.Range("F99").Value = 2054
$ws.Range("G99").Value = 1218
$ws.Range("H99").Value = 831
$ws.Range("I99").Value = 1188
$ws.Range("J99").Value = 992

# --- Washington Park ---
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("C4").Value = 5
$ws.Range("C6").Value = 22

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 12

# --- Uptown ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("D8").Value = 12
$ws.Range("E8").Value = 16
$ws.Range("D9").Value = 14
$ws.Range("E9").Value = 24

# --- Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("C7").Value = 33
$ws.Range("F8").Value = 43
$ws.Range("C9").Value = 97
$ws.Range("F9").Value = 82

# --- Lake View ---
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("F3").Value = 2
$ws.Range("H7").Value = 9
$ws.Range("J7").Value = 18
$ws.Range("F8").Value = 53
$ws.Range("H8").Value = 27
$ws.Range("J8").Value = 25

# --- Jefferson Park ---
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("E6").Value = 13
$ws.Range("E7").Value = 14

# --- Rogers Park ---
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("D2").Value = 2
$ws.Range("D7").Value = 31
$ws.Range("E7").Value = 54
$ws.Range("F7").Value = 29
$ws.Range("D8").Value = 41
$ws.Range("E8").Value = 68
$ws.Range("F8").Value = 44

# --- River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("B6").Value = 32
$ws.Range("B7").Value = 34

# --- West Loop ---
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I3").Value = 2
$ws.Range("I9").Value = 25

# --- Woodlawn ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("B2").Value = 2
$ws.Range("B7").Value = 21

# --- Near South Side ---
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 17
$ws.Range("F6").Value = 14
$ws.Range("D7").Value = 19
$ws.Range("F7").Value = 18

# --- Douglas ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 13

# --- West Pullman ---
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("F4").Value = 4
$ws.Range("B5").Value = 6
$ws.Range("B6").Value = 13
$ws.Range("F6").Value = 8

# --- Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("E3").Value = 5
$ws.Range("C9").Value = 29
$ws.Range("J9").Value = 25
$ws.Range("C10").Value = 46
$ws.Range("E10").Value = 52
$ws.Range("J10").Value = 50

# --- West Town ---
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("B5").Value = 2
$ws.Range("B7").Value = 10

# --- Lincoln Park ---
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("G6").Value = 6
$ws.Range("G7").Value = 7
$ws.Range("F3").Value = 2

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("B6").Value = 22
$ws.Range("H6").Value = 3
$ws.Range("B7").Value = 35
$ws.Range("F7").Value = 15
$ws.Range("H7").Value = 8

# --- O'Hare ---
$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("F5").Value = 26
$ws.Range("F6").Value = 28

# --- Gage Park ---
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("C7").Value = 18
$ws.Range("C8").Value = 21

# --- Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("E6").Value = 24
$ws.Range("C7").Value = 35
$ws.Range("F7").Value = 56
$ws.Range("C8").Value = 69
$ws.Range("E8").Value = 72
$ws.Range("F8").Value = 90

# --- Irving Park ---
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 3

# --- Kenwood ---
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J7").Value = 3

# --- Washington Heights ---
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("B5").Value = 1
$ws.Range("B7").Value = 6

# --- Oakland ---
$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 2
